$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (Test Type: add float to orchestration)
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "amount"
$ws.Range("D2").Value = "Ab tt"
$ws.Range("E2").Value = "INVALID_TYPE"
$ws.Range("F2").Value = "Value is not a valid float"

# Shorten regex failure messages across affected rows (17-19, 21-31)
$rows = @(17,18,19,21,22,23,24,25,26,27,28,29,30,31)
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "Regex failure"
}
